$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.002495288848877
$ws.Range("B1").Value = 2.112586498260498
$ws.Range("C1").Value = 6.793773651123047
$ws.Range("D1").Value = 1.952247977256775
$ws.Range("E1").Value = 1.372190833091736
